{"js": "// Update the date line and the 25 division problems in the practice table.\n// The edit is a pure text substitution: every run keeps its existing\n// run/paragraph formatting - only the literal text content changes.\n\nconst titleOld = \"2024-01-01 Monday\";\nconst titleNew = \"2024-01-02 Tuesday\";\n\n// Row-major (top-left -> bottom-right) replacement values for the 25\n// non-empty table cells (5 data rows of 5 cells each - the table also has\n// blank spacer rows in between which are left untouched).\nconst newValues = [\n  [\"71\u00f78=\", \"91\u00f77=\", \"56\u00f78=\", \"52\u00f76=\", \"86\u00f72=\"],\n  [\"88\u00f74=\", \"11\u00f76=\", \"86\u00f73=\", \"21\u00f79=\", \"88\u00f79=\"],\n  [\"71\u00f73=\", \"93\u00f74=\", \"33\u00f78=\", \"49\u00f77=\", \"90\u00f75=\"],\n  [\"45\u00f72=\", \"85\u00f72=\", \"61\u00f78=\", \"70\u00f79=\", \"94\u00f75=\"],\n  [\"92\u00f72=\", \"65\u00f72=\", \"78\u00f74=\", \"23\u00f77=\", \"98\u00f78=\"],\n];\n\nconst body = context.document.body;\n\n// First paragraph in the body holds the centered date heading.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nif (titlePara.text.trim() === titleOld) {\n  titlePara.getRange().insertText(titleNew, \"Replace\");\n}\n\n// Walk the (single) table and overwrite each data row's cells in order,\n// skipping the blank spacer rows.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rows = table.rows;\nrows.load(\"items/cellCount\");\nawait context.sync();\n\nlet dataRowIndex = 0;\nfor (let r = 0; r < rows.items.length && dataRowIndex < newValues.length; r++) {\n  const cellCount = rows.items[r].cellCount;\n  if (!cellCount) continue;\n\n  const firstCellBody = table.getCell(r, 0).body;\n  firstCellBody.load(\"text\");\n  await context.sync();\n\n  if (!firstCellBody.text || firstCellBody.text.trim() === \"\") continue;\n\n  const rowValues = newValues[dataRowIndex];\n  for (let c = 0; c < cellCount && c < rowValues.length; c++) {\n    table.getCell(r, c).body.getRange().insertText(rowValues[c], \"Replace\");\n  }\n  dataRowIndex++;\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the practice table.\n# The edit is a pure text substitution: every run keeps its existing\n# run/paragraph formatting - only the literal text content changes.\n\n$d = $word.ActiveDocument\n\n$titleOld = \"2024-01-01 Monday\"\n$titleNew = \"2024-01-02 Tuesday\"\n\n$titlePara = $d.Paragraphs.Item(1)\nif ($titlePara.Range.Text.TrimEnd(\"`r\") -eq $titleOld) {\n    $titlePara.Range.Text = $titleNew\n}\n\n# Row-major (top-left -> bottom-right) replacement values for the 25\n# non-empty table cells (5 rows of 5 cells each).\n$newValues = @(\n    @(\"71\u00f78=\", \"91\u00f77=\", \"56\u00f78=\", \"52\u00f76=\", \"86\u00f72=\"),\n    @(\"88\u00f74=\", \"11\u00f76=\", \"86\u00f73=\", \"21\u00f79=\", \"88\u00f79=\"),\n    @(\"71\u00f73=\", \"93\u00f74=\", \"33\u00f78=\", \"49\u00f77=\", \"90\u00f75=\"),\n    @(\"45\u00f72=\", \"85\u00f72=\", \"61\u00f78=\", \"70\u00f79=\", \"94\u00f75=\"),\n    @(\"92\u00f72=\", \"65\u00f72=\", \"78\u00f74=\", \"23\u00f77=\", \"98\u00f78=\")\n)\n\n$t = $d.Tables.Item(1)\n$dataRowIndex = 0\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $firstCellText = $t.Cell($r, 1).Range.Text.TrimEnd(\"`r\", [char]7)\n    if ($firstCellText.Trim() -eq \"\") {\n        continue\n    }\n    if ($dataRowIndex -ge $newValues.Count) {\n        break\n    }\n    $rowValues = $newValues[$dataRowIndex]\n    for ($c = 1; $c -le $t.Columns.Count -and $c -le $rowValues.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n    $dataRowIndex++\n}\n"}
